# Scheduled-runner update: refresh cached market-board figures (currentAveragePrice*,
# LevePrice*, LeveProfit*) for a handful of leve rows across each crafting-job sheet.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2050.6
$ws.Range("I33").Value = 793.1429000000001
$ws.Range("J33").Value = 4984.6665
$ws.Range("K33").Value = 793.1429000000001
$ws.Range("L33").Value = 4984.6665
$ws.Range("M33").Value = -564.1429000000001
$ws.Range("N33").Value = -5442.6665

$ws.Range("H40").Value = 3470.5454
$ws.Range("J40").Value = 1648.375
$ws.Range("L40").Value = 1648.375
$ws.Range("N40").Value = -1998.375

$ws.Range("H62").Value = 55560820
$ws.Range("I62").Value = 125001350
$ws.Range("J62").Value = 8397.200000000001
$ws.Range("K62").Value = 125001350
$ws.Range("L62").Value = 8397.200000000001
$ws.Range("M62").Value = -125000726
$ws.Range("N62").Value = -9645.200000000001

$ws.Range("H65").Value = 55560820
$ws.Range("I65").Value = 125001350
$ws.Range("J65").Value = 8397.200000000001
$ws.Range("K65").Value = 625006750
$ws.Range("L65").Value = 41986
$ws.Range("M65").Value = -625003630
$ws.Range("N65").Value = -48226

$ws.Range("H76").Value = 5400.4
$ws.Range("I76").Value = 4668.3335
$ws.Range("J76").Value = 6498.5
$ws.Range("K76").Value = 4668.3335
$ws.Range("L76").Value = 6498.5
$ws.Range("M76").Value = -4353.3335
$ws.Range("N76").Value = -7128.5

$ws.Range("H79").Value = 5400.4
$ws.Range("I79").Value = 4668.3335
$ws.Range("J79").Value = 6498.5
$ws.Range("K79").Value = 4668.3335
$ws.Range("L79").Value = 6498.5
$ws.Range("M79").Value = -3576.3335
$ws.Range("N79").Value = -8682.5

$ws.Range("H106").Value = 4799.5557
$ws.Range("I106").Value = 1742.5
$ws.Range("K106").Value = 1742.5
$ws.Range("M106").Value = -1111.5

$ws.Range("H132").Value = 5574.8066
$ws.Range("I132").Value = 4944.0415
$ws.Range("K132").Value = 14832.1245
$ws.Range("M132").Value = -12302.1245

$ws.Range("H138").Value = 3324.1553
$ws.Range("I138").Value = 1722.2307
$ws.Range("J138").Value = 3786.9333
$ws.Range("K138").Value = 5166.6921
$ws.Range("L138").Value = 11360.7999
$ws.Range("M138").Value = -26.69210000000021
$ws.Range("N138").Value = -21640.7999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1558
$ws.Range("I88").Value = 275.66666
$ws.Range("J88").Value = 2327.4
$ws.Range("K88").Value = 275.66666
$ws.Range("L88").Value = 2327.4
$ws.Range("M88").Value = 130.33334
$ws.Range("N88").Value = -3139.4

$ws.Range("H91").Value = 1558
$ws.Range("I91").Value = 275.66666
$ws.Range("J91").Value = 2327.4
$ws.Range("K91").Value = 275.66666
$ws.Range("L91").Value = 2327.4
$ws.Range("M91").Value = 1128.33334
$ws.Range("N91").Value = -5135.4

$ws.Range("H97").Value = 1022.2857
$ws.Range("I97").Value = 838.7059
$ws.Range("K97").Value = 838.7059
$ws.Range("M97").Value = -342.7059

$ws.Range("H110").Value = 1833.7391
$ws.Range("I110").Value = 1482.6666
$ws.Range("K110").Value = 1482.6666
$ws.Range("M110").Value = 562.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 7800
$ws.Range("I15").Value = 4333.3335
$ws.Range("K15").Value = 4333.3335
$ws.Range("M15").Value = -4106.3335

$ws.Range("H19").Value = 4769.6
$ws.Range("I19").Value = 3855.3333
$ws.Range("K19").Value = 3855.3333
$ws.Range("M19").Value = -3682.3333

$ws.Range("H80").Value = 501.10526
$ws.Range("I80").Value = 823.6667
$ws.Range("J80").Value = 352.23077
$ws.Range("K80").Value = 823.6667
$ws.Range("L80").Value = 352.23077
$ws.Range("M80").Value = 174.3333
$ws.Range("N80").Value = -2348.23077

$ws.Range("H82").Value = 39441.2
$ws.Range("I82").Value = 21068.666
$ws.Range("K82").Value = 21068.666
$ws.Range("M82").Value = -20685.666

$ws.Range("H83").Value = 501.10526
$ws.Range("I83").Value = 823.6667
$ws.Range("J83").Value = 352.23077
$ws.Range("K83").Value = 4118.3335
$ws.Range("L83").Value = 1761.15385
$ws.Range("M83").Value = 873.6665000000003
$ws.Range("N83").Value = -11745.15385

$ws.Range("H85").Value = 39441.2
$ws.Range("I85").Value = 21068.666
$ws.Range("K85").Value = 21068.666
$ws.Range("M85").Value = -19742.666

$ws.Range("H105").Value = 2807.625
$ws.Range("I105").Value = 2807.625
$ws.Range("K105").Value = 2807.625
$ws.Range("M105").Value = -1060.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 3502657
$ws.Range("I6").Value = 4903200
$ws.Range("K6").Value = 4903200
$ws.Range("M6").Value = -4903087

$ws.Range("H19").Value = 716.6667
$ws.Range("I19").Value = 750
$ws.Range("K19").Value = 750
$ws.Range("M19").Value = -580

$ws.Range("H24").Value = 716.6667
$ws.Range("I24").Value = 750
$ws.Range("K24").Value = 750
$ws.Range("M24").Value = -580

$ws.Range("H25").Value = 4980.778
$ws.Range("I25").Value = 2565.4
$ws.Range("K25").Value = 2565.4
$ws.Range("M25").Value = -2391.4

$ws.Range("H122").Value = 1588.5
$ws.Range("I122").Value = 1251.5
$ws.Range("J122").Value = 1925.5
$ws.Range("K122").Value = 3754.5
$ws.Range("L122").Value = 5776.5
$ws.Range("M122").Value = -1304.5
$ws.Range("N122").Value = -10676.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 8622094
$ws.Range("J131").Value = 1473.9608
$ws.Range("L131").Value = 4421.8824
$ws.Range("N131").Value = -14501.8824

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 14166.5
$ws.Range("I58").Value = 14166.5
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 14166.5
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -13889.5
$ws.Range("N58").ClearContents()

$ws.Range("H126").Value = 4947.9443
$ws.Range("I126").Value = 6277
$ws.Range("K126").Value = 18831
$ws.Range("M126").Value = -16361

$ws.Range("H132").Value = 80480.62
$ws.Range("I132").Value = 80480.62
$ws.Range("K132").Value = 241441.86
$ws.Range("M132").Value = -238911.86

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H31").Value = 1984.1428
$ws.Range("J31").Value = 2847.25
$ws.Range("L31").Value = 2847.25
$ws.Range("N31").Value = -3343.25

$ws.Range("H38").Value = 62799.6
$ws.Range("I38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("M38").ClearContents()

$ws.Range("H40").Value = 3367.0476
$ws.Range("I40").Value = 3076.5
$ws.Range("J40").Value = 4296.8
$ws.Range("K40").Value = 3076.5
$ws.Range("L40").Value = 4296.8
$ws.Range("M40").Value = -2940.5
$ws.Range("N40").Value = -4568.8

$ws.Range("H41").Value = 9999
$ws.Range("I41").Value = 9999
$ws.Range("K41").Value = 9999
$ws.Range("M41").Value = -9561

$ws.Range("H46").Value = 4137.2
$ws.Range("I46").Value = 3546.25
$ws.Range("J46").Value = 6501
$ws.Range("K46").Value = 3546.25
$ws.Range("L46").Value = 6501
$ws.Range("M46").Value = -3358.25
$ws.Range("N46").Value = -6877

$ws.Range("H50").Value = 38500
$ws.Range("I50").Value = 38000
$ws.Range("J50").Value = 39000
$ws.Range("K50").Value = 38000
$ws.Range("L50").Value = 39000
$ws.Range("M50").Value = -37363
$ws.Range("N50").Value = -40274

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 10015.5
$ws.Range("I34").Value = 10015.5
$ws.Range("K34").Value = 10015.5
$ws.Range("M34").Value = -9812.5

$ws.Range("H37").Value = 74943
$ws.Range("I37").Value = 49998
$ws.Range("J37").Value = 99888
$ws.Range("K37").Value = 49998
$ws.Range("L37").Value = 99888
$ws.Range("M37").Value = -49795
$ws.Range("N37").Value = -100294

$ws.Range("H40").Value = 47499
$ws.Range("I40").Value = 47499
$ws.Range("K40").Value = 47499
$ws.Range("M40").Value = -47350

$ws.Range("H42").Value = 49998
$ws.Range("I42").Value = 49998
$ws.Range("K42").Value = 49998
$ws.Range("M42").Value = -49620

$ws.Range("H43").Value = 39999
$ws.Range("I43").Value = 39999
$ws.Range("K43").Value = 39999
$ws.Range("M43").Value = -39850

$ws.Range("H49").Value = 49998
$ws.Range("I49").Value = 49998
$ws.Range("K49").Value = 49998
$ws.Range("M49").Value = -49768

$ws.Range("H122").Value = 13566.148
$ws.Range("I122").Value = 3333.25
$ws.Range("K122").Value = 9999.75
$ws.Range("M122").Value = -7549.75

$ws.Range("H126").Value = 2689.7727
$ws.Range("I126").Value = 2765.1765
$ws.Range("J126").Value = 2433.4
$ws.Range("K126").Value = 8295.529500000001
$ws.Range("L126").Value = 7300.200000000001
$ws.Range("M126").Value = -5825.529500000001
$ws.Range("N126").Value = -12240.2
